# Nexial showcase workbook update:
#   [base] - [`outputToCloud(resource)`]: support the transferring of output
#   artifact to the cloud.
#
# This adds a brand-new "text" category (with its single function
# `spellCheck(var,profile,text)`) to the '#system' reference sheet, and adds
# a new `outputToCloud(resource)` function to the existing "base" category.
#
# Concretely, on the hidden '#system' sheet:
#   1. A new column is inserted at Y for the "text" category, which pushes
#      the existing Y..AD columns (web, webalert, webcookie, ws, ws.async,
#      xml) one column to the right (Z..AE).
#   2. A new row is inserted into column A (the "target"/category list) so
#      that "text" takes its alphabetical spot just before "web".
#   3. A new row is inserted into column E (the "base" function list) so
#      that `outputToCloud(resource)` takes its alphabetical spot just
#      before `prependText(...)`.
#   4. The workbook-level defined names are updated/added to match the new
#      ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$lastDataRow = 129

# ---------------------------------------------------------------------------
# Step 1: make room for the new "text" column by shifting columns
# Y(25)..AD(30) one column to the right, into Z(26)..AE(31). Walk from the
# right-most column to the left-most so a column's data is safely copied out
# before it gets overwritten by its left neighbour.
# ---------------------------------------------------------------------------
for ($srcCol = 30; $srcCol -ge 25; $srcCol--) {
    $dstCol = $srcCol + 1

    $colValues = @()
    for ($r = 1; $r -le $lastDataRow; $r++) {
        $colValues += $ws.Cells.Item($r, $srcCol).Value()
    }

    for ($r = 1; $r -le $lastDataRow; $r++) {
        $ws.Cells.Item($r, $dstCol).Value = ""
    }

    for ($r = 1; $r -le $lastDataRow; $r++) {
        $v = $colValues[$r - 1]
        if ($v -ne $null) {
            $ws.Cells.Item($r, $dstCol).Value = $v
        }
    }
}

# Column Y (25) is now the brand-new "text" category: a header plus its one
# function.
for ($r = 1; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, 25).Value = ""
}
$ws.Cells.Item(1, 25).Value = "text"
$ws.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# Step 2: insert "text" into the category list (column A) right before
# "web", shifting web/webalert/webcookie/ws/ws.async/xml down by one row.
# Walk bottom-up so each value is read before being overwritten.
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------------------
# Step 3: insert `outputToCloud(resource)` into the "base" function list
# (column E) right before `prependText(...)`, shifting the remaining
# functions down by one row.
# ---------------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $v = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r + 1, 5).Value = $v
}
$ws.Cells.Item(22, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# Step 4: update the workbook-level defined names so their ranges match the
# rows/columns that just moved, and register the brand-new "text" name.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
